$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 456; this shifts the existing rows 456-556
# down to 457-557 (matching the dimension change A1:R556 -> A1:R557).
$ws.Rows.Item(456).Insert()

# Populate the newly inserted row 456 with the new record.
$ws.Range("A456").Value = 8
$ws.Range("B456").Value = "Terminal La Palmera de La Serena"
$ws.Range("C456").Value = "Coquimbo"
$ws.Range("D456").Value = 44641
$ws.Range("E456").Value = 4
$ws.Range("F456").Value = 100112024
$ws.Range("G456").Value = "Choclo"
$ws.Range("H456").Value = "Dulce o Americano"
$ws.Range("I456").Value = "Primera"
$ws.Range("J456").Value = 460
$ws.Range("K456").Value = 13500
$ws.Range("L456").Value = 14000
$ws.Range("M456").Value = 13750
$ws.Range("N456").Value = "`$/malla 70 unidades"
$ws.Range("O456").Value = "Provincia del Elquí"
$ws.Range("P456").Value = 196
$ws.Range("Q456").Value = 70
$ws.Range("R456").Value = "Hortaliza"
